# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Replace the worker/period/value rows (B16:G24) with the updated data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: DocType, DocNumber, Name, Period, ValorMora, SalarioBasico
$data = @(
    @("CC", "1128058778", "YOSIMAR LORDUY CHAVEZ", "1706", 29509, 737717),
    @("CC", "1128058778", "YOSIMAR LORDUY CHAVEZ", "1705", 29509, 737717),
    @("CC", "1128058778", "YOSIMAR LORDUY CHAVEZ", "1704", 29509, 737717),
    @("CC", "1128058778", "YOSIMAR LORDUY CHAVEZ", "1703", 29509, 737717),
    @("CC", "73188736",   "OSCAR EUGENIO ALVAREZ ANGULO",    "1608", 100454, 781242),
    @("CC", "9294622",    "OSVALDO ENRIQUE ALVAREZ MARTINEZ","1701", 96000,  781242),
    @("CC", "73270962",   "EDEL ENRIQUE CHAMORRO CANTILLO",  "1610", 32000,  781242),
    @("CC", "73270962",   "EDEL ENRIQUE CHAMORRO CANTILLO",  "1609", 19200,  781242),
    @("CC", "1049482102", "YOHON DAIRO DE LEON BERROCAL",    "1609", 27734,  781242)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
}
